# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. "总计" summary sheet: insert a new row 2 for "2022-Q4", push the
#    existing "2022-Q3" row down to row 3.
# ----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Row 2 ("A2", label style bold+centered+border) formatting also belongs
# on the relocated row 3 -- copy it down before writing values there.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0

# ----------------------------------------------------------------------
# 2. Add a brand-new "2022-Q4" worksheet positioned right before
#    "2022-Q3" (so tab order becomes 总计, 2022-Q4, 2022-Q3).
# ----------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# Reuse the same header/label formatting (bold, centered, bordered) that
# the "总计" sheet already uses for its header row and its "A" column.
$wsTotal.Range("B1:D1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)

# Header row
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# The B:G data columns hold numeric-looking figures that must stay text
# (e.g. "0.30", "0.0022") so their original formatting/precision survives
# -- format them as Text before assigning so Excel doesn't re-parse them
# as numbers. (G4's source value is a genuine 0, so it is left out here
# and kept a real number below.)
$wsQ4.Range("B2:G3").NumberFormat = "@"
$wsQ4.Range("B4:F4").NumberFormat = "@"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").Value = "850007"
$wsQ4.Range("C2").Value = "海通智选一年持有期股票B"
$wsQ4.Range("D2").Value = "0.30"
$wsQ4.Range("E2").Value = "82.33"
$wsQ4.Range("F2").Value = "0.73"
$wsQ4.Range("G2").Value = "0.0022"
$wsQ4.Range("H2").Value = 6

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").Value = "850788"
$wsQ4.Range("C3").Value = "海通智选一年持有期股票A"
$wsQ4.Range("D3").Value = "0.20"
$wsQ4.Range("E3").Value = "82.33"
$wsQ4.Range("F3").Value = "0.73"
$wsQ4.Range("G3").Value = "0.0015"
$wsQ4.Range("H3").Value = 6

$wsQ4.Range("A4").Value = 2
$wsQ4.Range("B4").Value = "850799"
$wsQ4.Range("C4").Value = "海通智选一年持有期股票C"
$wsQ4.Range("D4").Value = "0.00"
$wsQ4.Range("E4").Value = "82.33"
$wsQ4.Range("F4").Value = "0.73"

# G4's underlying value is exactly 0, so (matching the source data) it is
# left as a genuine number instead of a formatted "0.00"-style string.
$wsQ4.Range("G4").Value = 0
$wsQ4.Range("H4").Value = 6

# "2022-Q3" kept the workbook's active/selected tab before this edit --
# re-resolve it by name (its index shifted when "2022-Q4" was inserted
# ahead of it) and restore it as the active sheet.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Activate()
